$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.038.09"
$ws.Range("E2").Value = "  -0.03%  "
$ws.Range("D3").Value = "2.299.84"
$ws.Range("E3").Value = "  -0.01%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'300.21"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.14%  "
$ws.Range("D6").Value = "'97.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.69%  "
$ws.Range("D7").Value = "'0.520"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.54%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  +0.61%  "
$ws.Range("D10").Value = "'36.07"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.72%  "
$ws.Range("D11").Value = "'0.0790"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("E12").Value = "  +0.76%  "
$ws.Range("D13").Value = "'17.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.91%  "
$ws.Range("E14").Value = "  -0.55%  "
$ws.Range("D15").Value = "2.658.20"
$ws.Range("E15").Value = "  +0.01%  "
$ws.Range("D16").Value = "2.324.31"
$ws.Range("E16").Value = "  +2.04%  "
$ws.Range("E17").Value = "  -1.44%  "
$ws.Range("D18").Value = "42.933.17"
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "'12.75"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "0.0₃0908"
$ws.Range("E20").Value = "  +0.26%  "
$ws.Range("D21").Value = "'6.13"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Value = "'68.74"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.15%  "
$ws.Range("D23").Value = "'237.85"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.86%  "
$ws.Range("E24").Value = "  -1.00%  "
$ws.Range("E25").Value = "  -0.35%  "
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("E27").Value = "  -0.41%  "
$ws.Range("E28").Value = "  +0.18%  "
$ws.Range("D29").Value = "'164.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.67%  "
$ws.Range("E31").Value = "  -0.11%  "
$ws.Range("D32").Value = "'33.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.28%  "
$ws.Range("E33").Value = "  +0.02%  "
$ws.Range("E34").Value = "  +1.26%  "
$ws.Range("E35").Value = "  +4.00%  "
$ws.Range("D36").Value = "'17.95"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.71%  "
$ws.Range("E37").Value = "  +0.10%  "
$ws.Range("D38").Value = "'0.0698"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.21%  "
$ws.Range("E39").Value = "  +0.56%  "
$ws.Range("E40").Value = "  -0.57%  "
$ws.Range("E41").Value = "  -1.29%  "
$ws.Range("E42").Value = "  +0.91%  "
$ws.Range("D43").Value = "2.015.45"
$ws.Range("E43").Value = "  +1.52%  "
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("E45").Value = "  -3.96%  "
$ws.Range("D46").Value = "'10.40"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.25%  "
$ws.Range("D47").Value = "'17.45"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.61%  "
$ws.Range("E48").Value = "  -2.50%  "
$ws.Range("E49").Value = "  -2.97%  "
$ws.Range("D50").Value = "2.528.32"
$ws.Range("E50").Value = "  +0.16%  "
$ws.Range("E51").Value = "  -1.45%  "
